$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new posting row (row 5): reuse the existing "ML Engineer/Data Scientist"
# job description text from A2, with Min/Max years experience of 2 and 4.
$ws.Range("A5").Value2 = $ws.Range("A2").Value2
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 4

# Re-fit the row height so it matches the default (undo any auto row-height bump
# triggered by inserting multi-line text into a brand-new row).
$ws.Rows.Item(5).AutoFit()
